$wb = $excel.ActiveWorkbook

# --- Step 1: insert new summary row in "总计" sheet for 2022-Q3 ---
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A2").Value = 0
$summary.Range('A3').Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 19
$summary.Range("D2").Value = 4.91

# --- Step 2: create the new "2022-Q3" sheet by duplicating "2022-Q2" (keeps formatting) ---
$firstSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $firstSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# --- Step 3: overwrite the fund-holding rows with 2022-Q3 data ---
$q3Sheet.Range("B2").Value = "'000362"
$q3Sheet.Range("B2").ClearFormats()
$q3Sheet.Range("C2").Value = "国泰聚信价值优势灵活配置混合A"
$q3Sheet.Range("D2").Value = "'27.52"
$q3Sheet.Range("D2").ClearFormats()
$q3Sheet.Range("E2").Value = "'89.04"
$q3Sheet.Range("E2").ClearFormats()
$q3Sheet.Range("F2").Value = "'4.80"
$q3Sheet.Range("F2").ClearFormats()
$q3Sheet.Range("G2").Value = "'1.3210"
$q3Sheet.Range("G2").ClearFormats()
$q3Sheet.Range("H2").Value = 4

$q3Sheet.Range("B3").Value = "'000363"
$q3Sheet.Range("B3").ClearFormats()
$q3Sheet.Range("C3").Value = "国泰聚信价值优势灵活配置混合C"
$q3Sheet.Range("D3").Value = "'13.07"
$q3Sheet.Range("D3").ClearFormats()
$q3Sheet.Range("E3").Value = "'89.04"
$q3Sheet.Range("E3").ClearFormats()
$q3Sheet.Range("F3").Value = "'4.80"
$q3Sheet.Range("F3").ClearFormats()
$q3Sheet.Range("G3").Value = "'0.6274"
$q3Sheet.Range("G3").ClearFormats()
$q3Sheet.Range("H3").Value = 4

$q3Sheet.Range("B4").Value = "'001579"
$q3Sheet.Range("B4").ClearFormats()
$q3Sheet.Range("C4").Value = "国泰大农业股票A"
$q3Sheet.Range("D4").Value = "'12.15"
$q3Sheet.Range("D4").ClearFormats()
$q3Sheet.Range("E4").Value = "'88.79"
$q3Sheet.Range("E4").ClearFormats()
$q3Sheet.Range("F4").Value = "'4.62"
$q3Sheet.Range("F4").ClearFormats()
$q3Sheet.Range("G4").Value = "'0.5613"
$q3Sheet.Range("G4").ClearFormats()
$q3Sheet.Range("H4").Value = 4

$q3Sheet.Range("B5").Value = "'008415"
$q3Sheet.Range("B5").ClearFormats()
$q3Sheet.Range("C5").Value = "国泰大制造两年持有期混合"
$q3Sheet.Range("D5").Value = "'10.30"
$q3Sheet.Range("D5").ClearFormats()
$q3Sheet.Range("E5").Value = "'90.83"
$q3Sheet.Range("E5").ClearFormats()
$q3Sheet.Range("F5").Value = "'4.96"
$q3Sheet.Range("F5").ClearFormats()
$q3Sheet.Range("G5").Value = "'0.5109"
$q3Sheet.Range("G5").ClearFormats()
$q3Sheet.Range("H5").Value = 3

$q3Sheet.Range("B6").Value = "'020010"
$q3Sheet.Range("B6").ClearFormats()
$q3Sheet.Range("C6").Value = "国泰金牛创新混合"
$q3Sheet.Range("D6").Value = "'13.26"
$q3Sheet.Range("D6").ClearFormats()
$q3Sheet.Range("E6").Value = "'86.27"
$q3Sheet.Range("E6").ClearFormats()
$q3Sheet.Range("F6").Value = "'2.94"
$q3Sheet.Range("F6").ClearFormats()
$q3Sheet.Range("G6").Value = "'0.3898"
$q3Sheet.Range("G6").ClearFormats()
$q3Sheet.Range("H6").Value = 9

$q3Sheet.Range("B7").Value = "'012173"
$q3Sheet.Range("B7").ClearFormats()
$q3Sheet.Range("C7").Value = "国泰兴泽优选一年持有期混合A"
$q3Sheet.Range("D7").Value = "'8.41"
$q3Sheet.Range("D7").ClearFormats()
$q3Sheet.Range("E7").Value = "'88.23"
$q3Sheet.Range("E7").ClearFormats()
$q3Sheet.Range("F7").Value = "'4.45"
$q3Sheet.Range("F7").ClearFormats()
$q3Sheet.Range("G7").Value = "'0.3742"
$q3Sheet.Range("G7").ClearFormats()
$q3Sheet.Range("H7").Value = 3

$q3Sheet.Range("B8").Value = "'012174"
$q3Sheet.Range("B8").ClearFormats()
$q3Sheet.Range("C8").Value = "国泰兴泽优选一年持有期混合C"
$q3Sheet.Range("D8").Value = "'6.17"
$q3Sheet.Range("D8").ClearFormats()
$q3Sheet.Range("E8").Value = "'88.23"
$q3Sheet.Range("E8").ClearFormats()
$q3Sheet.Range("F8").Value = "'4.45"
$q3Sheet.Range("F8").ClearFormats()
$q3Sheet.Range("G8").Value = "'0.2746"
$q3Sheet.Range("G8").ClearFormats()
$q3Sheet.Range("H8").Value = 3

$q3Sheet.Range("B9").Value = "'007835"
$q3Sheet.Range("B9").ClearFormats()
$q3Sheet.Range("C9").Value = "国泰鑫睿混合"
$q3Sheet.Range("D9").Value = "'8.30"
$q3Sheet.Range("D9").ClearFormats()
$q3Sheet.Range("E9").Value = "'79.49"
$q3Sheet.Range("E9").ClearFormats()
$q3Sheet.Range("F9").Value = "'3.13"
$q3Sheet.Range("F9").ClearFormats()
$q3Sheet.Range("G9").Value = "'0.2598"
$q3Sheet.Range("G9").ClearFormats()
$q3Sheet.Range("H9").Value = 9

$q3Sheet.Range("B10").Value = "'013890"
$q3Sheet.Range("B10").ClearFormats()
$q3Sheet.Range("C10").Value = "国泰睿毅三年持有期混合A"
$q3Sheet.Range("D10").Value = "'4.86"
$q3Sheet.Range("D10").ClearFormats()
$q3Sheet.Range("E10").Value = "'89.26"
$q3Sheet.Range("E10").ClearFormats()
$q3Sheet.Range("F10").Value = "'3.86"
$q3Sheet.Range("F10").ClearFormats()
$q3Sheet.Range("G10").Value = "'0.1876"
$q3Sheet.Range("G10").ClearFormats()
$q3Sheet.Range("H10").Value = 9

$q3Sheet.Range("B11").Value = "'005244"
$q3Sheet.Range("B11").ClearFormats()
$q3Sheet.Range("C11").Value = "国泰聚优价值灵活配置混合A"
$q3Sheet.Range("D11").Value = "'4.61"
$q3Sheet.Range("D11").ClearFormats()
$q3Sheet.Range("E11").Value = "'87.30"
$q3Sheet.Range("E11").ClearFormats()
$q3Sheet.Range("F11").Value = "'4.00"
$q3Sheet.Range("F11").ClearFormats()
$q3Sheet.Range("G11").Value = "'0.1844"
$q3Sheet.Range("G11").ClearFormats()
$q3Sheet.Range("H11").Value = 5

$q3Sheet.Range("B12").Value = "'005245"
$q3Sheet.Range("B12").ClearFormats()
$q3Sheet.Range("C12").Value = "国泰聚优价值灵活配置混合C"
$q3Sheet.Range("D12").Value = "'1.80"
$q3Sheet.Range("D12").ClearFormats()
$q3Sheet.Range("E12").Value = "'87.30"
$q3Sheet.Range("E12").ClearFormats()
$q3Sheet.Range("F12").Value = "'4.00"
$q3Sheet.Range("F12").ClearFormats()
$q3Sheet.Range("G12").Value = "'0.0720"
$q3Sheet.Range("G12").ClearFormats()
$q3Sheet.Range("H12").Value = 5

$q3Sheet.Range("B13").Value = "'160324"
$q3Sheet.Range("B13").ClearFormats()
$q3Sheet.Range("C13").Value = "华夏磐晟灵活配置混合（LOF）"
$q3Sheet.Range("D13").Value = "'0.91"
$q3Sheet.Range("D13").ClearFormats()
$q3Sheet.Range("E13").Value = "'90.76"
$q3Sheet.Range("E13").ClearFormats()
$q3Sheet.Range("F13").Value = "'4.24"
$q3Sheet.Range("F13").ClearFormats()
$q3Sheet.Range("G13").Value = "'0.0386"
$q3Sheet.Range("G13").ClearFormats()
$q3Sheet.Range("H13").Value = 10

$q3Sheet.Range("B14").Value = "'000573"
$q3Sheet.Range("B14").ClearFormats()
$q3Sheet.Range("C14").Value = "天弘通利混合"
$q3Sheet.Range("D14").Value = "'1.01"
$q3Sheet.Range("D14").ClearFormats()
$q3Sheet.Range("E14").Value = "'79.25"
$q3Sheet.Range("E14").ClearFormats()
$q3Sheet.Range("F14").Value = "'3.45"
$q3Sheet.Range("F14").ClearFormats()
$q3Sheet.Range("G14").Value = "'0.0348"
$q3Sheet.Range("G14").ClearFormats()
$q3Sheet.Range("H14").Value = 4

$q3Sheet.Range("B15").Value = "'163110"
$q3Sheet.Range("B15").ClearFormats()
$q3Sheet.Range("C15").Value = "申万菱信量化小盘股票（LOF）A"
$q3Sheet.Range("D15").Value = "'5.04"
$q3Sheet.Range("D15").ClearFormats()
$q3Sheet.Range("E15").Value = "'93.06"
$q3Sheet.Range("E15").ClearFormats()
$q3Sheet.Range("F15").Value = "'0.59"
$q3Sheet.Range("F15").ClearFormats()
$q3Sheet.Range("G15").Value = "'0.0297"
$q3Sheet.Range("G15").ClearFormats()
$q3Sheet.Range("H15").Value = 10

$q3Sheet.Range("B16").Value = "'013891"
$q3Sheet.Range("B16").ClearFormats()
$q3Sheet.Range("C16").Value = "国泰睿毅三年持有期混合C"
$q3Sheet.Range("D16").Value = "'0.45"
$q3Sheet.Range("D16").ClearFormats()
$q3Sheet.Range("E16").Value = "'89.26"
$q3Sheet.Range("E16").ClearFormats()
$q3Sheet.Range("F16").Value = "'3.86"
$q3Sheet.Range("F16").ClearFormats()
$q3Sheet.Range("G16").Value = "'0.0174"
$q3Sheet.Range("G16").ClearFormats()
$q3Sheet.Range("H16").Value = 9

$q3Sheet.Range("B17").Value = "'004194"
$q3Sheet.Range("B17").ClearFormats()
$q3Sheet.Range("C17").Value = "招商中证1000指数增强A"
$q3Sheet.Range("D17").Value = "'1.56"
$q3Sheet.Range("D17").ClearFormats()
$q3Sheet.Range("E17").Value = "'92.06"
$q3Sheet.Range("E17").ClearFormats()
$q3Sheet.Range("F17").Value = "'1.03"
$q3Sheet.Range("F17").ClearFormats()
$q3Sheet.Range("G17").Value = "'0.0161"
$q3Sheet.Range("G17").ClearFormats()
$q3Sheet.Range("H17").Value = 10

$q3Sheet.Range("B18").Value = "'004195"
$q3Sheet.Range("B18").ClearFormats()
$q3Sheet.Range("C18").Value = "招商中证1000指数增强C"
$q3Sheet.Range("D18").Value = "'1.09"
$q3Sheet.Range("D18").ClearFormats()
$q3Sheet.Range("E18").Value = "'92.06"
$q3Sheet.Range("E18").ClearFormats()
$q3Sheet.Range("F18").Value = "'1.03"
$q3Sheet.Range("F18").ClearFormats()
$q3Sheet.Range("G18").Value = "'0.0112"
$q3Sheet.Range("G18").ClearFormats()
$q3Sheet.Range("H18").Value = 10

$q3Sheet.Range("B19").Value = "'015588"
$q3Sheet.Range("B19").ClearFormats()
$q3Sheet.Range("C19").Value = "国泰大农业股票C"
$q3Sheet.Range("D19").Value = "'0.07"
$q3Sheet.Range("D19").ClearFormats()
$q3Sheet.Range("E19").Value = "'88.79"
$q3Sheet.Range("E19").ClearFormats()
$q3Sheet.Range("F19").Value = "'4.62"
$q3Sheet.Range("F19").ClearFormats()
$q3Sheet.Range("G19").Value = "'0.0032"
$q3Sheet.Range("G19").ClearFormats()
$q3Sheet.Range("H19").Value = 4

$q3Sheet.Range("B20").Value = "'013918"
$q3Sheet.Range("B20").ClearFormats()
$q3Sheet.Range("C20").Value = "申万菱信量化小盘股票（LOF）C"
$q3Sheet.Range("D20").Value = "'0.00"
$q3Sheet.Range("D20").ClearFormats()
$q3Sheet.Range("E20").Value = "'93.06"
$q3Sheet.Range("E20").ClearFormats()
$q3Sheet.Range("F20").Value = "'0.59"
$q3Sheet.Range("F20").ClearFormats()
$q3Sheet.Range("G20").Value = 0
$q3Sheet.Range("H20").Value = 10

Write-Host "done"
